$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from an existing header cell
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF)
$dataI = @(9, 9, 8, 9, 6, 8, 8, 7, 6, 4, 8, 8, 5, 7)
$dataJ = @(9, 9, 9, 9, 6, 8, 8, 7, 6, 4, 8, 8, 5, 7)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
